$wb = $excel.ActiveWorkbook

# --- Sheet "Games": append the completed game that used to be the first
#     row of "Next" (OKC on 45306) as new row 42 ---
$games = $wb.Worksheets.Item("Games")
$r = 42
$games.Cells.Item($r,1).Value = 41
$games.Cells.Item($r,2).Value = 45306
$games.Cells.Item($r,2).NumberFormat = "YYYY-MM-DD"
$games.Cells.Item($r,3).Value = 1
$games.Cells.Item($r,4).Value = 112
$games.Cells.Item($r,5).Value = 97
$games.Cells.Item($r,6).Value = 0.571
$games.Cells.Item($r,7).Value = 15.6
$games.Cells.Item($r,8).Value = 22.9
$games.Cells.Item($r,9).Value = 0.19
$games.Cells.Item($r,10).Value = 115.5
$games.Cells.Item($r,11).Value = "OKC"
$games.Cells.Item($r,12).Value = 105
$games.Cells.Item($r,13).Value = 0.495
$games.Cells.Item($r,14).Value = 9.800000000000001
$games.Cells.Item($r,15).Value = 26.5
$games.Cells.Item($r,16).Value = 0.104
$games.Cells.Item($r,17).Value = 108.3
$games.Cells.Item($r,18).Value = 1
$games.Cells.Item($r,19).Value = 1

# --- Sheet "Next": that game is no longer upcoming, so drop its row and
#     shift the remaining upcoming games up by one ---
$next = $wb.Worksheets.Item("Next")
$next.Rows.Item(2).Delete()
